$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.648.05'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '3.394.12'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('D5').Value = '''560.63'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = '''175.65'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '''0.630'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').Value = '3.384.86'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('E10').Value = '  +3.93%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '''53.68'
$ws.Range('E12').Value = '  -2.35%  '
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').Value = '3.936.35'
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').Value = '''18.25'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '3.386.98'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '65.633.90'
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').Value = '''11.85'
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('D21').Value = '''0.999'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '''482.07'
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('D23').Value = '''4.97'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').Value = '''89.99'
$ws.Range('E24').Value = '  +3.57%  '
$ws.Range('E25').Value = '  +3.77%  '
$ws.Range('D26').Value = '''4.11'
$ws.Range('E26').Value = '  -1.42%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  -3.22%  '
$ws.Range('D29').Value = '''8.73'
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').Value = '''31.33'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').Value = '''6.57'
$ws.Range('E31').Value = '  -2.53%  '
$ws.Range('D32').Value = '''63.73'
$ws.Range('E32').Value = '  +5.47%  '
$ws.Range('D33').Value = '''11.44'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').Value = '''572.77'
$ws.Range('E34').Value = '  -2.47%  '
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '''3.65'
$ws.Range('E37').Value = '  +4.08%  '
$ws.Range('D38').Value = '''0.141'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').Value = '''35.80'
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('D40').Value = '''0.375'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').Value = '0.0₃0743'
$ws.Range('E41').Value = '  -1.74%  '
$ws.Range('D42').Value = '3.096.53'
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('D43').Value = '''2.81'
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('E46').Value = '  -3.44%  '
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').Value = '''140.50'
$ws.Range('E49').Value = '  +2.42%  '
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('E51').Value = '  +0.54%  '
